$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 11189.667
$ws.Cells.Item(32, 10).Value = 9932.066000000001
$ws.Cells.Item(32, 12).Value = 9932.066000000001
$ws.Cells.Item(32, 14).Value = -10584.066
$ws.Cells.Item(69, 8).Value = 7353.125
$ws.Cells.Item(69, 9).Value = 1997.5
$ws.Cells.Item(69, 10).Value = 9138.333000000001
$ws.Cells.Item(69, 11).Value = 5992.5
$ws.Cells.Item(69, 12).Value = 27414.999
$ws.Cells.Item(69, 13).Value = -5118.5
$ws.Cells.Item(69, 14).Value = -29162.999
$ws.Cells.Item(72, 8).Value = 7353.125
$ws.Cells.Item(72, 9).Value = 1997.5
$ws.Cells.Item(72, 10).Value = 9138.333000000001
$ws.Cells.Item(72, 11).Value = 17977.5
$ws.Cells.Item(72, 12).Value = 82244.997
$ws.Cells.Item(72, 13).Value = -13609.5
$ws.Cells.Item(72, 14).Value = -90980.997
$ws.Cells.Item(105, 8).Value = 41399
$ws.Cells.Item(105, 10).Value = 41399
$ws.Cells.Item(105, 12).Value = 41399
$ws.Cells.Item(105, 14).Value = -48387
$ws.Cells.Item(135, 8).Value = 1536.6666
$ws.Cells.Item(135, 9).Value = 1234.7693
$ws.Cells.Item(135, 11).Value = 11112.9237
$ws.Cells.Item(135, 13).Value = -8577.923699999999
$ws.Cells.Item(137, 8).Value = 1407.4482
$ws.Cells.Item(137, 9).Value = 1338.9615
$ws.Cells.Item(137, 10).Value = 2001
$ws.Cells.Item(137, 11).Value = 4016.8845
$ws.Cells.Item(137, 12).Value = 6003
$ws.Cells.Item(137, 13).Value = -1466.8845
$ws.Cells.Item(137, 14).Value = -11103
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4557.732
$ws.Cells.Item(32, 9).Value = 4123.4326
$ws.Cells.Item(32, 10).Value = 8575
$ws.Cells.Item(32, 11).Value = 4123.4326
$ws.Cells.Item(32, 12).Value = 8575
$ws.Cells.Item(32, 13).Value = -3836.4326
$ws.Cells.Item(32, 14).Value = -9149
$ws.Cells.Item(39, 8).Value = 3750
$ws.Cells.Item(39, 9).Value = 3750
$ws.Cells.Item(39, 11).Value = 3750
$ws.Cells.Item(39, 13).Value = -3230
$ws.Cells.Item(40, 8).Value = 10000
$ws.Cells.Item(40, 10).Value = 10000
$ws.Cells.Item(40, 12).Value = 10000
$ws.Cells.Item(40, 14).Value = -10352
$ws.Cells.Item(45, 8).Value = 9074.706
$ws.Cells.Item(45, 9).Value = 14847.25
$ws.Cells.Item(45, 11).Value = 14847.25
$ws.Cells.Item(45, 13).Value = -14470.25
$ws.Cells.Item(61, 8).Value = 3020.516
$ws.Cells.Item(61, 9).Value = 2698.889
$ws.Cells.Item(61, 10).Value = 5191.5
$ws.Cells.Item(61, 11).Value = 2698.889
$ws.Cells.Item(61, 12).Value = 5191.5
$ws.Cells.Item(61, 13).Value = -2486.889
$ws.Cells.Item(61, 14).Value = -5615.5
$ws.Cells.Item(110, 8).Value = 9003.4375
$ws.Cells.Item(110, 9).Value = 12119.588
$ws.Cells.Item(110, 10).Value = 5471.8
$ws.Cells.Item(110, 11).Value = 12119.588
$ws.Cells.Item(110, 12).Value = 5471.8
$ws.Cells.Item(110, 13).Value = -10074.588
$ws.Cells.Item(110, 14).Value = -9561.799999999999
$ws.Cells.Item(132, 8).Value = 3884.56
$ws.Cells.Item(132, 9).Value = 3334.15
$ws.Cells.Item(132, 11).Value = 10002.45
$ws.Cells.Item(132, 13).Value = -7472.450000000001
$ws.Cells.Item(134, 8).Value = 64633.332
$ws.Cells.Item(134, 10).Value = 64633.332
$ws.Cells.Item(134, 12).Value = 64633.332
$ws.Cells.Item(134, 14).Value = -74773.33199999999
$ws.Cells.Item(136, 8).Value = 3020.516
$ws.Cells.Item(136, 9).Value = 2698.889
$ws.Cells.Item(136, 10).Value = 5191.5
$ws.Cells.Item(136, 11).Value = 8096.667
$ws.Cells.Item(136, 12).Value = 15574.5
$ws.Cells.Item(136, 13).Value = -5546.667
$ws.Cells.Item(136, 14).Value = -20674.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(36, 8).Value = 15292.571
$ws.Cells.Item(36, 9).Value = 1174.6666
$ws.Cells.Item(36, 11).Value = 1174.6666
$ws.Cells.Item(36, 13).Value = -640.6666
$ws.Cells.Item(64, 8).Value = 1292.8667
$ws.Cells.Item(64, 9).Value = 1663.6666
$ws.Cells.Item(64, 10).Value = 1045.6666
$ws.Cells.Item(64, 11).Value = 1663.6666
$ws.Cells.Item(64, 12).Value = 1045.6666
$ws.Cells.Item(64, 13).Value = -1438.6666
$ws.Cells.Item(64, 14).Value = -1495.6666
$ws.Cells.Item(67, 8).Value = 1292.8667
$ws.Cells.Item(67, 9).Value = 1663.6666
$ws.Cells.Item(67, 10).Value = 1045.6666
$ws.Cells.Item(67, 11).Value = 1663.6666
$ws.Cells.Item(67, 12).Value = 1045.6666
$ws.Cells.Item(67, 13).Value = -883.6666
$ws.Cells.Item(67, 14).Value = -2605.6666
$ws.Cells.Item(94, 8).Value = 1815.6428
$ws.Cells.Item(94, 9).Value = 1539.5416
$ws.Cells.Item(94, 10).Value = 3472.25
$ws.Cells.Item(94, 11).Value = 1539.5416
$ws.Cells.Item(94, 12).Value = 3472.25
$ws.Cells.Item(94, 13).Value = -1088.5416
$ws.Cells.Item(94, 14).Value = -4374.25
$ws.Cells.Item(99, 8).Value = 5777.9
$ws.Cells.Item(99, 9).Value = 1800.4
$ws.Cells.Item(99, 11).Value = 1800.4
$ws.Cells.Item(99, 13).Value = -302.4000000000001
$ws.Cells.Item(103, 8).Value = 25413
$ws.Cells.Item(103, 10).Value = 25413
$ws.Cells.Item(103, 12).Value = 25413
$ws.Cells.Item(103, 14).Value = -27757
$ws.Cells.Item(132, 8).Value = 82500
$ws.Cells.Item(132, 10).Value = 82500
$ws.Cells.Item(132, 12).Value = 82500
$ws.Cells.Item(132, 14).Value = -92620
$ws.Cells.Item(134, 8).Value = 2466.4666
$ws.Cells.Item(134, 9).Value = 2160.3157
$ws.Cells.Item(134, 11).Value = 6480.9471
$ws.Cells.Item(134, 13).Value = -3945.9471
$ws.Cells.Item(139, 8).Value = 57025.6
$ws.Cells.Item(139, 9).Value = 43000
$ws.Cells.Item(139, 10).Value = 60532
$ws.Cells.Item(139, 11).Value = 43000
$ws.Cells.Item(139, 12).Value = 60532
$ws.Cells.Item(139, 13).Value = -37860
$ws.Cells.Item(139, 14).Value = -70812
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(32, 8).Value = 8040.8
$ws.Cells.Item(32, 9).Value = 12702
$ws.Cells.Item(32, 10).Value = 3379.6
$ws.Cells.Item(32, 11).Value = 12702
$ws.Cells.Item(32, 12).Value = 3379.6
$ws.Cells.Item(32, 13).Value = -12386
$ws.Cells.Item(32, 14).Value = -4011.6
$ws.Cells.Item(35, 8).Value = 690
$ws.Cells.Item(35, 9).Value = 483.33334
$ws.Cells.Item(35, 10).Value = 1000
$ws.Cells.Item(35, 11).Value = 483.33334
$ws.Cells.Item(35, 12).Value = 1000
$ws.Cells.Item(35, 13).Value = -189.33334
$ws.Cells.Item(35, 14).Value = -1588
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 60000
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 60000
$ws.Cells.Item(68, 13).Value = $null
$ws.Cells.Item(68, 14).Value = -61498
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 60000
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 180000
$ws.Cells.Item(71, 13).Value = $null
$ws.Cells.Item(71, 14).Value = -187488
$ws.Cells.Item(132, 8).Value = 4620.7334
$ws.Cells.Item(132, 9).Value = 4258.7
$ws.Cells.Item(132, 11).Value = 12776.1
$ws.Cells.Item(132, 13).Value = -10246.1
$ws.Cells.Item(137, 8).Value = 89055.44500000001
$ws.Cells.Item(137, 9).Value = 85000
$ws.Cells.Item(137, 10).Value = 89562.375
$ws.Cells.Item(137, 11).Value = 85000
$ws.Cells.Item(137, 12).Value = 89562.375
$ws.Cells.Item(137, 13).Value = -79900
$ws.Cells.Item(137, 14).Value = -99762.375
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(109, 8).Value = 486.5
$ws.Cells.Item(109, 9).Value = 486.5
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 11).Value = 1459.5
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 13).Value = -419.5
$ws.Cells.Item(109, 14).Value = $null
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(18, 8).Value = 13904398
$ws.Cells.Item(18, 9).Value = 18521798
$ws.Cells.Item(18, 11).Value = 18521798
$ws.Cells.Item(18, 13).Value = -18521505
$ws.Cells.Item(21, 8).Value = 4750001.5
$ws.Cells.Item(21, 9).Value = 4750001.5
$ws.Cells.Item(21, 11).Value = 4750001.5
$ws.Cells.Item(21, 13).Value = -4749828.5
$ws.Cells.Item(30, 8).Value = 4750001.5
$ws.Cells.Item(30, 9).Value = 4750001.5
$ws.Cells.Item(30, 11).Value = 4750001.5
$ws.Cells.Item(30, 13).Value = -4749896.5
$ws.Cells.Item(113, 8).Value = 1738.4445
$ws.Cells.Item(113, 9).Value = 1738.4445
$ws.Cells.Item(113, 11).Value = 1738.4445
$ws.Cells.Item(113, 13).Value = 431.5554999999999
$ws.Cells.Item(122, 8).Value = 2493.2307
$ws.Cells.Item(122, 9).Value = 2338.4443
$ws.Cells.Item(122, 10).Value = 2841.5
$ws.Cells.Item(122, 11).Value = 7015.3329
$ws.Cells.Item(122, 12).Value = 8524.5
$ws.Cells.Item(122, 13).Value = -4565.3329
$ws.Cells.Item(122, 14).Value = -13424.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1161.1111
$ws.Cells.Item(22, 10).Value = 1332.1428
$ws.Cells.Item(22, 12).Value = 1332.1428
$ws.Cells.Item(22, 14).Value = -1922.1428
$ws.Cells.Item(23, 8).Value = 8505998
$ws.Cells.Item(23, 9).Value = 503000
$ws.Cells.Item(23, 10).Value = 12507498
$ws.Cells.Item(23, 11).Value = 503000
$ws.Cells.Item(23, 12).Value = 12507498
$ws.Cells.Item(23, 13).Value = -502770
$ws.Cells.Item(23, 14).Value = -12507958
$ws.Cells.Item(27, 8).Value = 1161.1111
$ws.Cells.Item(27, 10).Value = 1332.1428
$ws.Cells.Item(27, 12).Value = 1332.1428
$ws.Cells.Item(27, 14).Value = -1546.1428
$ws.Cells.Item(32, 8).Value = 4413.875
$ws.Cells.Item(32, 9).Value = 2187.2856
$ws.Cells.Item(32, 10).Value = 20000
$ws.Cells.Item(32, 11).Value = 2187.2856
$ws.Cells.Item(32, 12).Value = 20000
$ws.Cells.Item(32, 13).Value = -1870.2856
$ws.Cells.Item(32, 14).Value = -20634
$ws.Cells.Item(35, 8).Value = 5938.6665
$ws.Cells.Item(35, 9).Value = 2854
$ws.Cells.Item(35, 10).Value = 9023.333000000001
$ws.Cells.Item(35, 11).Value = 2854
$ws.Cells.Item(35, 12).Value = 9023.333000000001
$ws.Cells.Item(35, 13).Value = -2518
$ws.Cells.Item(35, 14).Value = -9695.333000000001
$ws.Cells.Item(46, 8).Value = 1724
$ws.Cells.Item(46, 9).Value = 1598.5
$ws.Cells.Item(46, 10).Value = 1975
$ws.Cells.Item(46, 11).Value = 1598.5
$ws.Cells.Item(46, 12).Value = 1975
$ws.Cells.Item(46, 13).Value = -1410.5
$ws.Cells.Item(46, 14).Value = -2351
$ws.Cells.Item(61, 8).Value = 3673.9524
$ws.Cells.Item(61, 10).Value = 3600.8
$ws.Cells.Item(61, 12).Value = 3600.8
$ws.Cells.Item(61, 14).Value = -4004.8
$ws.Cells.Item(113, 8).Value = 3673.9524
$ws.Cells.Item(113, 10).Value = 3600.8
$ws.Cells.Item(113, 12).Value = 3600.8
$ws.Cells.Item(113, 14).Value = -7940.8
$ws.Cells.Item(132, 8).Value = 4888.706
$ws.Cells.Item(132, 9).Value = 3773
$ws.Cells.Item(132, 11).Value = 11319
$ws.Cells.Item(132, 13).Value = -8789
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(30, 8).Value = 12290.643
$ws.Cells.Item(30, 10).Value = 12290.643
$ws.Cells.Item(30, 12).Value = 12290.643
$ws.Cells.Item(30, 14).Value = -12504.643
$ws.Cells.Item(42, 8).Value = 28133.334
$ws.Cells.Item(42, 10).Value = 28133.334
$ws.Cells.Item(42, 12).Value = 28133.334
$ws.Cells.Item(42, 14).Value = -28889.334
$ws.Cells.Item(100, 8).Value = 1845.2354
$ws.Cells.Item(100, 9).Value = 1918
$ws.Cells.Item(100, 10).Value = 1794.3
$ws.Cells.Item(100, 11).Value = 3836
$ws.Cells.Item(100, 12).Value = 3588.6
$ws.Cells.Item(100, 13).Value = -3295
$ws.Cells.Item(100, 14).Value = -4670.6
$ws.Cells.Item(132, 8).Value = 3515.6875
$ws.Cells.Item(132, 9).Value = 3487.6785
$ws.Cells.Item(132, 11).Value = 10463.0355
$ws.Cells.Item(132, 13).Value = -7933.0355
